$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.666.81"
$ws.Range("E2").Value = "  -4.86%  "

$ws.Range("D3").Value = "2.985.81"
$ws.Range("E3").Value = "  -5.71%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.63"
$ws.Range("E5").Value = "  -4.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.41"
$ws.Range("E6").Value = "  -8.08%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  -4.84%  "

$ws.Range("D9").Value = "2.971.61"
$ws.Range("E9").Value = "  -6.09%  "

$ws.Range("E10").Value = "  -5.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.20"
$ws.Range("E11").Value = "  -7.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.364"
$ws.Range("E12").Value = "  -5.60%  "

$ws.Range("D13").Value = "3.501.53"
$ws.Range("E13").Value = "  -5.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.123"
$ws.Range("E14").Value = "  -3.71%  "

$ws.Range("D15").Value = "61.778.26"
$ws.Range("E15").Value = "  -4.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.54"
$ws.Range("E16").Value = "  -7.62%  "

$ws.Range("D17").Value = "2.987.09"
$ws.Range("E17").Value = "  -5.59%  "

$ws.Range("E18").Value = "  -5.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "389.92"
$ws.Range("E19").Value = "  -5.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.08"
$ws.Range("E20").Value = "  -4.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.88"
$ws.Range("E21").Value = "  -6.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.59"
$ws.Range("E22").Value = "  -7.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.69"
$ws.Range("E24").Value = "  -5.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.464"
$ws.Range("E25").Value = "  -4.59%  "

$ws.Range("E26").Value = "  -8.89%  "

$ws.Range("E27").Value = "  +0.26%  "

$ws.Range("D28").Value = "0.0₃0937"
$ws.Range("E28").Value = "  -9.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.40"
$ws.Range("E29").Value = "  -5.30%  "

$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.72"
$ws.Range("E31").Value = "  -5.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.32"
$ws.Range("E32").Value = "  -4.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.49"
$ws.Range("E33").Value = "  +2.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.63"
$ws.Range("E34").Value = "  -5.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.95"
$ws.Range("E35").Value = "  -6.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.07"
$ws.Range("E36").Value = "  -5.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.27"
$ws.Range("E37").Value = "  -6.06%  "

$ws.Range("E38").Value = "  -9.13%  "

$ws.Range("D39").Value = "2.434.25"
$ws.Range("E39").Value = "  -10.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.43"
$ws.Range("E40").Value = "  -3.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.89"
$ws.Range("E41").Value = "  -5.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.17"
$ws.Range("E42").Value = "  -6.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.656"
$ws.Range("E43").Value = "  -7.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0591"
$ws.Range("E44").Value = "  -6.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.996"
$ws.Range("E45").Value = "  -0.38%  "

$ws.Range("E46").Value = "  -6.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.87"
$ws.Range("E47").Value = "  -12.18%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.68"
$ws.Range("E48").Value = "  -7.70%  "

$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.47"
$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0942"
$ws.Range("E50").Value = "  -4.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "262.96"
$ws.Range("E51").Value = "  -9.81%  "
